# Generate Report for Handback
# Refresh the handoff / handback timestamps recorded for the
# "771d07f9-2832-47ca-9120-017e14e31e6f.md" row (row 3) across the
# Overview / zh-cn / de-de report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2017-02-21 08:45:49"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (L) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2017-02-21 08:45:31"
$zhcn.Range("L3").Value = "2017-02-21 08:47:11"

# --- de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (L) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2017-02-21 08:45:49"
$dede.Range("L3").Value = "2017-02-21 08:47:34"
